$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "63.220.31"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +4.05%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.055.23"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.06%  "
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "550.25"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +4.29%  "
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "138.74"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +5.90%  "
$ws.Range("E7").Value = "  -0.02%  "
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "3.047.90"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("E11").Value = "  +0.58%  "
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.450"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +2.18%  "
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000225"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +3.32%  "
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "34.56"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +3.37%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.551.70"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +2.53%  "
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "63.252.46"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +3.90%  "
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "3.054.77"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  +3.02%  "
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "479.41"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +4.59%  "
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "13.55"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +3.49%  "
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "0.672"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "7.15"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +4.78%  "
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "80.92"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +3.22%  "
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "12.39"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +3.89%  "
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "7.84"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +6.76%  "
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "25.89"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +2.75%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "1.15"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  +6.96%  "
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "5.64"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +5.91%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "55.54"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "5.95"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "461.28"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  +4.11%  "
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "3.111.35"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("E41").Value = "  +0.51%  "
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "8.20"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +2.03%  "
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "2.53"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +4.62%  "
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "27.86"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +7.56%  "
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.250"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("E46").Value = "  -0.12%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "2.03"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("E48").Value = "  +1.51%  "
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "115.78"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -2.46%  "
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0505"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("E51").Value = "  +4.79%  "
